$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the raw "count" values in column B that feed the formulas in
# columns D/F. These were previously blank (formulas evaluated to 0).
$ws.Range("B181").Value = 1998
$ws.Range("B182").Value = 6
$ws.Range("B183").Value = 530
$ws.Range("B184").Value = 470
$ws.Range("B185").Value = 352
$ws.Range("B189").Value = 484
$ws.Range("B190").Value = 78

# Force a full recalculation so dependent formulas (D181, F181, ... F191)
# pick up the new values.
$excel.CalculateFullRebuild()

# Update the view so it matches where the user was scrolled to / selected
# (top-left visible cell A175, cursor on J199).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 175
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J199").Select()
